$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely; remaining columns (B:F) shift left to (A:E)
$ws.Range("A:A").Delete()
